$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.934.14'
$ws.Range("E2").Value = '  -0.28%  '

$ws.Range("D3").Value = '1.823.44'
$ws.Range("E3").Value = '  +0.03%  '

$ws.Range("E4").Value = '  -0.69%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.17'
$ws.Range("E5").Value = '  +0.15%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  -0.63%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4626'
$ws.Range("E7").Value = '  +0.17%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3688'
$ws.Range("E8").Value = '  +1.28%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07323'
$ws.Range("E9").Value = '  +0.53%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8743'
$ws.Range("E10").Value = '  +1.08%  '

$ws.Range("E11").Value = '  +3.29%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '19.64'
$ws.Range("E12").Value = '  -1.01%  '

$ws.Range("D13").Value = '1.809.66'
$ws.Range("E13").Value = '  -1.39%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.558'
$ws.Range("E14").Value = '  +0.99%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.331'
$ws.Range("E15").Value = '  +0.00%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.16'
$ws.Range("E16").Value = '  -2.25%  '

$ws.Range("E17").Value = '  -0.35%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008834'
$ws.Range("E18").Value = '  +2.41%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.003'
$ws.Range("E19").Value = '  -0.56%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.73'
$ws.Range("E20").Value = '  +1.64%  '

$ws.Range("D21").Value = '26.964.22'
$ws.Range("E21").Value = '  -1.60%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.097'
$ws.Range("E22").Value = '  -1.30%  '

$ws.Range("E23").Value = '  -0.66%  '

$ws.Range("D24").Value = '2.064.93'
$ws.Range("E24").Value = '  -2.17%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.80'
$ws.Range("E25").Value = '  +0.11%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.851'
$ws.Range("E26").Value = '  -0.26%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.33'
$ws.Range("E27").Value = '  +0.67%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.034'
$ws.Range("E28").Value = '  -2.67%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.096'
$ws.Range("E29").Value = '  -0.06%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.36'
$ws.Range("E30").Value = '  -0.46%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08867'
$ws.Range("E31").Value = '  -0.37%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.956'
$ws.Range("E32").Value = '  +0.20%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7307'
$ws.Range("E33").Value = '  +0.37%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.428'
$ws.Range("E34").Value = '  +0.00%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.130'
$ws.Range("E35").Value = '  -0.95%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.470'
$ws.Range("E36").Value = '  -1.26%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.072'
$ws.Range("E37").Value = '  -0.04%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01937'
$ws.Range("E38").Value = '  +0.99%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05208'
$ws.Range("E39").Value = '  -1.22%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.954'
$ws.Range("E40").Value = '  +0.93%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.085'
$ws.Range("E41").Value = '  -1.20%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5129'
$ws.Range("E42").Value = '  -1.48%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1620'
$ws.Range("E43").Value = '  -0.85%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.138'

$ws.Range("E45").Value = '  -0.80%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.20'
$ws.Range("E46").Value = '  +0.85%  '

$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.003'
$ws.Range("E47").Value = '  -0.68%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '101.75'
$ws.Range("E48").Value = '  -1.35%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.618'
$ws.Range("E49").Value = '  -0.85%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06196'
$ws.Range("E50").Value = '  -0.49%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '64.47'
$ws.Range("E51").Value = '  -0.27%  '
